# Commit: "added ability to have custom suffix replacement"
#
# The workbook's "data" sheet stores harmonized variable names in column D
# (rows 2-17) built from a "...|Harmonized-DB" suffix. This edit renames
# that suffix to "...|Harmonized" (dropping the "-DB"), reflecting that the
# suffix used when harmonizing is now configurable rather than hard-coded.
#
# Reproduce this the way a user would in Excel: select the column of
# variable names on the "data" sheet and run Find & Replace (Ctrl+H) for
# "Harmonized-DB" -> "Harmonized".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data")
$ws.Activate()

$rng = $ws.Range("D2:D17")

$rng.Replace("Harmonized-DB", "Harmonized")

# Leave the just-edited range selected, as it would be right after running
# Find & Replace on that selection.
$rng.Select()
